$p = $ppt.ActivePresentation

# --- Slide 3 (sldId 263 / cId 167852364): "Heart attacks: they're bad" ---
# Extend the 2nd bullet's sentence about cardiovascular-disease deaths.
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(2)
$tr3 = $sh3.TextFrame.TextRange
$full3 = $tr3.Text
$needle3 = "Major cardiovascular disease is recorded as causing over 600,000 deaths annually, most of which are due to heart attacks."
$idx3 = $full3.IndexOf($needle3)
$c3 = $tr3.Characters($idx3 + 1, $needle3.Length)
$c3.Text = "Major cardiovascular disease is recorded as causing over 600,000 deaths annually, most of which are due to heart attacks, making it the CVDs responsible for one in four deaths in America."

# --- Slide 6 (sldId 260 / cId 2548556833): "What's new" ---
# Finish the sentence describing what the clinical data input includes.
$s6 = $p.Slides.Item(6)
$sh6 = $s6.Shapes.Item(2)
$tr6 = $sh6.TextFrame.TextRange
$tr6.Text = "The purpose of my research is to incorporate electrocardiogram test results into the input of the predictive algorithm along with the clinical data in an effort to increase prediction accuracy. Clinical data input would include a person’s age, cholesterol levels, gender, history of smoking, body mass index, and diabetes. Input would also include a number of cycles of an electrocardiogram."

# --- Slide 7 (sldId 261 / cId 278273717): "Procedure" ---
# Nudge the rotated right-arrow connector (id 12) to its new position.
$s7 = $p.Slides.Item(7)
$sh7 = $s7.Shapes.Item(7)
$sh7.Left = 678.4936
$sh7.Top = 195.4047
